$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 2")

for ($col = 2; $col -le 29; $col++) {
    $cell1 = $ws.Cells.Item(127, $col)
    $cell2 = $ws.Cells.Item(128, $col)
    $val1 = $cell1.Value2
    $val2 = $cell2.Value2
    $cell1.Value2 = $val2
    $cell2.Value2 = $val1
}

for ($col = 2; $col -le 29; $col++) {
    $cell1 = $ws.Cells.Item(130, $col)
    $cell2 = $ws.Cells.Item(131, $col)
    $val1 = $cell1.Value2
    $val2 = $cell2.Value2
    $cell1.Value2 = $val2
    $cell2.Value2 = $val1
}
